# Append two new log rows (98 and 99) to the trading log worksheet,
# matching the "Update trading results" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 98: TRADING_ATTEMPT for ETH
$ws.Cells.Item(98, 1).Value = "2025-10-28T01:28:07.082034"
$ws.Cells.Item(98, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(98, 3).Value = "ETH"
$ws.Cells.Item(98, 4).Value = "UNKNOWN"
$ws.Cells.Item(98, 5).Value = 4128.421039015593
$ws.Cells.Item(98, 11).Value = "ATTEMPT"
$ws.Cells.Item(98, 12).Value = "Attempting trade 1/1"

# Row 99: POSITION_FAILED for ETH
$ws.Cells.Item(99, 1).Value = "2025-10-28T01:28:08.970846"
$ws.Cells.Item(99, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(99, 3).Value = "ETH"
$ws.Cells.Item(99, 4).Value = "UNKNOWN"
$ws.Cells.Item(99, 11).Value = "FAILED"
$ws.Cells.Item(99, 12).Value = "Trade execution failed for trade 1"
